$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-07 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-08 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("481×9=", $true, $false, $false, $false, $false, $true, 1, $false, "451×3=", 2) | Out-Null
$d.Content.Find.Execute("436×2=", $true, $false, $false, $false, $false, $true, 1, $false, "438×9=", 2) | Out-Null
$d.Content.Find.Execute("665×6=", $true, $false, $false, $false, $false, $true, 1, $false, "596×4=", 2) | Out-Null
$d.Content.Find.Execute("315×7=", $true, $false, $false, $false, $false, $true, 1, $false, "222×3=", 2) | Out-Null
$d.Content.Find.Execute("296×5=", $true, $false, $false, $false, $false, $true, 1, $false, "680×4=", 2) | Out-Null
$d.Content.Find.Execute("858×3=", $true, $false, $false, $false, $false, $true, 1, $false, "956×5=", 2) | Out-Null
$d.Content.Find.Execute("940×3=", $true, $false, $false, $false, $false, $true, 1, $false, "114×8=", 2) | Out-Null
$d.Content.Find.Execute("481×8=", $true, $false, $false, $false, $false, $true, 1, $false, "918×8=", 2) | Out-Null
$d.Content.Find.Execute("716×5=", $true, $false, $false, $false, $false, $true, 1, $false, "107×7=", 2) | Out-Null
$d.Content.Find.Execute("704×9=", $true, $false, $false, $false, $false, $true, 1, $false, "520×3=", 2) | Out-Null
$d.Content.Find.Execute("649×4=", $true, $false, $false, $false, $false, $true, 1, $false, "795×7=", 2) | Out-Null
$d.Content.Find.Execute("946×7=", $true, $false, $false, $false, $false, $true, 1, $false, "192×6=", 2) | Out-Null
$d.Content.Find.Execute("301×9=", $true, $false, $false, $false, $false, $true, 1, $false, "850×4=", 2) | Out-Null
$d.Content.Find.Execute("695×2=", $true, $false, $false, $false, $false, $true, 1, $false, "200×6=", 2) | Out-Null
$d.Content.Find.Execute("194×5=", $true, $false, $false, $false, $false, $true, 1, $false, "441×5=", 2) | Out-Null
$d.Content.Find.Execute("506×2=", $true, $false, $false, $false, $false, $true, 1, $false, "732×9=", 2) | Out-Null
$d.Content.Find.Execute("588×3=", $true, $false, $false, $false, $false, $true, 1, $false, "444×4=", 2) | Out-Null
$d.Content.Find.Execute("140×3=", $true, $false, $false, $false, $false, $true, 1, $false, "542×6=", 2) | Out-Null
$d.Content.Find.Execute("405×9=", $true, $false, $false, $false, $false, $true, 1, $false, "773×5=", 2) | Out-Null
$d.Content.Find.Execute("108×7=", $true, $false, $false, $false, $false, $true, 1, $false, "149×3=", 2) | Out-Null
$d.Content.Find.Execute("986×2=", $true, $false, $false, $false, $false, $true, 1, $false, "167×9=", 2) | Out-Null
$d.Content.Find.Execute("240×6=", $true, $false, $false, $false, $false, $true, 1, $false, "171×2=", 2) | Out-Null
$d.Content.Find.Execute("562×3=", $true, $false, $false, $false, $false, $true, 1, $false, "518×5=", 2) | Out-Null
$d.Content.Find.Execute("481×7=", $true, $false, $false, $false, $false, $true, 1, $false, "902×2=", 2) | Out-Null
$d.Content.Find.Execute("716×8=", $true, $false, $false, $false, $false, $true, 1, $false, "157×6=", 2) | Out-Null

$d.Save()
